# Refresh cryptos list values (price/volume columns), per GitHub Actions update.
# Row 34/35 additionally swap Coin name + Link because the ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (prices, percentages)
    # are not reinterpreted as numbers - matches the source inlineStr cells.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" '36.431.71'
Set-TextCell "E2" '  +0.40%  '
Set-TextCell "D3" '1.943.38'
Set-TextCell "E3" '  -1.13%  '
Set-TextCell "E4" '  -0.07%  '
Set-TextCell "D5" '243.13'
Set-TextCell "E5" '  -0.42%  '
Set-TextCell "E6" '  -1.44%  '
Set-TextCell "E7" '  -0.11%  '
Set-TextCell "D8" '57.44'
Set-TextCell "E8" '  -0.25%  '
Set-TextCell "D9" '0.361'
Set-TextCell "E9" '  -2.48%  '
Set-TextCell "D10" '0.0845'
Set-TextCell "E10" '  -1.35%  '
Set-TextCell "D11" '0.102'
Set-TextCell "E11" '  -1.53%  '
Set-TextCell "D12" '2.228.44'
Set-TextCell "E12" '  -1.19%  '
Set-TextCell "D13" '21.43'
Set-TextCell "E13" '  -2.37%  '
Set-TextCell "D14" '0.813'
Set-TextCell "E14" '  -3.07%  '
Set-TextCell "D15" '13.50'
Set-TextCell "E15" '  -0.42%  '
Set-TextCell "D16" '5.16'
Set-TextCell "E16" '  -3.51%  '
Set-TextCell "D17" '1.952.24'
Set-TextCell "E17" '  -0.52%  '
Set-TextCell "D18" '36.382.44'
Set-TextCell "E18" '  +0.61%  '
Set-TextCell "D19" '69.27'
Set-TextCell "E19" '  -2.47%  '
Set-TextCell "D20" '0.0₃0863'
Set-TextCell "E20" '  -2.69%  '
Set-TextCell "D21" '228.66'
Set-TextCell "E21" '  -1.74%  '
Set-TextCell "D22" '5.01'
Set-TextCell "E22" '  -2.91%  '
Set-TextCell "E23" '  -0.23%  '
Set-TextCell "D24" '2.36'
Set-TextCell "E24" '  -5.68%  '
Set-TextCell "D25" '2.30'
Set-TextCell "E25" '  +0.83%  '
Set-TextCell "D26" '9.23'
Set-TextCell "E26" '  -4.13%  '
Set-TextCell "D27" '161.95'
Set-TextCell "E27" '  -2.11%  '
Set-TextCell "E28" '  +3.24%  '
Set-TextCell "D29" '19.22'
Set-TextCell "E29" '  -3.72%  '
Set-TextCell "E30" '  -0.80%  '
Set-TextCell "D31" '1.10'
Set-TextCell "E31" '  -5.36%  '
Set-TextCell "D32" '4.58'
Set-TextCell "E32" '  -4.27%  '
Set-TextCell "D33" '0.0618'
Set-TextCell "E33" '  -3.85%  '
Set-TextCell "B34" 'THORChain'
Set-TextCell "C34" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell "D34" '6.26'
Set-TextCell "E34" '  +5.05%  '
Set-TextCell "B35" 'InternetComputer(DFINITY)'
Set-TextCell "C35" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell "D35" '4.19'
Set-TextCell "E35" '  -4.19%  '
Set-TextCell "E36" '  +0.13%  '
Set-TextCell "E37" '  -0.73%  '
Set-TextCell "D38" '2.17'
Set-TextCell "E38" '  +1.03%  '
Set-TextCell "D39" '3.14'
Set-TextCell "E39" '  +6.92%  '
Set-TextCell "D40" '0.0989'
Set-TextCell "E40" '  +3.28%  '
Set-TextCell "E41" '  +0.00%  '
Set-TextCell "E42" '  -0.91%  '
Set-TextCell "D43" '1.15'
Set-TextCell "E43" '  -3.28%  '
Set-TextCell "D44" '15.95'
Set-TextCell "E44" '  +1.16%  '
Set-TextCell "D45" '1.342.45'
Set-TextCell "E45" '  -0.21%  '
Set-TextCell "D46" '1.03'
Set-TextCell "E46" '  -3.56%  '
Set-TextCell "D47" '86.80'
Set-TextCell "E47" '  -2.68%  '
Set-TextCell "D48" '7.17'
Set-TextCell "E48" '  -1.85%  '
Set-TextCell "E49" '  +0.29%  '
Set-TextCell "D50" '2.119.38'
Set-TextCell "E50" '  -1.10%  '
Set-TextCell "D51" '43.29'
Set-TextCell "E51" '  -2.99%  '
